$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112 (shifts rows 112:140 down to 113:141)
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with the new weekly record
$ws.Cells.Item(112, 1).Value = 8
$ws.Cells.Item(112, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(112, 3).Value = "Coquimbo"
$ws.Cells.Item(112, 4).Value = 45093
$ws.Cells.Item(112, 5).Value = 4
$ws.Cells.Item(112, 6).Value = 100114007
$ws.Cells.Item(112, 7).Value = "Jengibre"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 360
$ws.Cells.Item(112, 11).Value = 17000
$ws.Cells.Item(112, 12).Value = 18000
$ws.Cells.Item(112, 13).Value = 17500
$ws.Cells.Item(112, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(112, 15).Value = "Perú"
$ws.Cells.Item(112, 16).Value = 1346
$ws.Cells.Item(112, 17).Value = 13
$ws.Cells.Item(112, 18).Value = "Hortaliza"
